$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Tracker row 9: fill in the two newly-entered data values.
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 6

# Move/save the cursor selection to G24, matching the latest save state.
$ws.Activate()
$ws.Range("G24").Select()
